$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.579.02'
$ws.Range('E2').Value = '  -1.36%  '
$ws.Range('D3').Value = '3.363.07'
$ws.Range('E3').Value = '  -2.76%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '555.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '175.95'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.620'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.54%  '
$ws.Range('D8').Value = '3.350.59'
$ws.Range('E8').Value = '  -2.92%  '
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.631'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.164'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.68'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000273'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.07'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.37%  '
$ws.Range('D15').Value = '3.888.04'
$ws.Range('E15').Value = '  -3.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.42'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.118'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.41%  '
$ws.Range('D18').Value = '3.362.83'
$ws.Range('E18').Value = '  -2.75%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.87'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.75%  '
$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D20').Value = '64.461.50'
$ws.Range('E20').Value = '  -1.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.984'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '454.14'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +10.93%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.91'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +11.70%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.08'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.40'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.41'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.80'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.84'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.79'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.80%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '29.97'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.64'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.52'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '578.47'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.108'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '58.51'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.93%  '
$ws.Range('E36').Value = '  +0.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.141'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -8.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.50'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.67'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.03%  '
$ws.Range('D40').Value = '0.0₃0755'
$ws.Range('E40').Value = '  -4.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.368'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.22%  '
$ws.Range('D42').Value = '3.101.62'
$ws.Range('E42').Value = '  -3.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.998'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.79'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.34%  '
$ws.Range('E45').Value = '  -2.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0410'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.48'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.131'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.58'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.33'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '135.45'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.94%  '
